$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1 (matching the style of the existing H1 header)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill columns I and J for data rows 2 through 38
# Column I is a constant flag of 1, column J mirrors column H's value
for ($r = 2; $r -le 38; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $hVal = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 10).Value = $hVal
}
